$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark two more items as reserved (column E = "Reserved")
$ws.Range("E3").Value = "Y"
$ws.Range("E17").Value = "Y"

# Turn the plain-text product URL in C17 into a real hyperlink,
# matching the styling used by the other links in column C.
$ws.Hyperlinks.Add($ws.Range("C17"), "https://www.amazon.de/dp/B06XHQ3G3C?psc=1&ref_=cm_sw_r_cp_ud_ct_WRC2J9TPX505376KF519")
$ws.Range("C17").Style = "Hyperlink"

# Move the active selection
[void]$ws.Range("E8").Select()
